{"js": "// Merge the split \"improving transparency\" / \"and adding modern structure\" / \".\"\n// runs back into a single run of text, and fix the AMEC date range from\n// \"08/2013 - 04/2014\" to \"03/2013 - 03/2014\".\n\nconst body = context.document.body;\n\n// --- Change 1: merge the three runs about \"improving transparency ...\" ---\nconst searchResults1 = body.search(\"improving transparency and adding modern structure.\", { matchCase: true });\nsearchResults1.load(\"items\");\nawait context.sync();\n\nif (searchResults1.items.length === 0) {\n  // Fallback: the text might still be split across runs (unlikely after search,\n  // since Word's search matches visible text regardless of run boundaries).\n  throw new Error(\"Could not find target text for change 1\");\n}\n\nsearchResults1.items[0].insertText(\n  \"improving transparency and adding modern structure.\",\n  Word.InsertLocation.replace\n);\n\n// --- Change 2: fix the AMEC Earth & Environmental date range ---\nconst searchResults2 = body.search(\"08/2013 - 04/2014\", { matchCase: true });\nsearchResults2.load(\"items\");\nawait context.sync();\n\nif (searchResults2.items.length === 0) {\n  throw new Error(\"Could not find target text for change 2\");\n}\n\nsearchResults2.items[0].insertText(\"03/2013 - 03/2014\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Merge the split \"improving transparency\" / \"and adding modern structure\" / \".\"\n# runs back into a single run of text, and fix the AMEC date range from\n# \"08/2013 - 04/2014\" to \"03/2013 - 03/2014\".\n\n$d = $word.ActiveDocument\n\n# --- Change 1: merge the three runs about \"improving transparency ...\" ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"Reducing labor costs, improving transparency and adding modern structure.\"\n$find1.Replacement.Text = \"Reducing labor costs, improving transparency and adding modern structure.\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# --- Change 2: fix the AMEC Earth & Environmental date range ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"08/2013 - 04/2014\"\n$find2.Replacement.Text = \"03/2013 - 03/2014\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
